# Generate Report for Archive
#
# The underlying report generator re-walked the e2e folder and this time
# encountered 431a5951-24a4-4fd0-a3e7-a2afb768907c.md before
# a9819c07-1654-4020-80cc-c104e8a674a9.md, so those two files' rows
# (row 3 and row 4) trade places on every sheet. The other rows
# (320f1920..., row 2, and 6e85a030..., row 5) are unaffected.
#
# For every sheet we only touch the handful of cells that actually carry
# per-file data (file name / path columns, status, source-xliff name and
# its datetime) - leaving shared/style-identical columns (Source Path,
# Priority, Content Duplicate, booleans, blanks, etc.) completely alone so
# their cell types/shared-string slots aren't disturbed.
#
# The hyperlinks for column B (Overview) / column A (zh-cn, de-de) keep
# pointing at the same r:id (i.e. the same external GitHub URL) for a
# given row, but their displayed text swaps along with the cell text.

function Swap-RowCells {
    param(
        $ws,
        [int]$row1,
        [int]$row2,
        [string[]]$cols
    )
    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

function Swap-HyperlinkDisplay {
    param(
        $ws,
        [int]$row1,
        [int]$row2
    )
    $text1 = $null
    $text2 = $null
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq $row1) { $text1 = $hl.TextToDisplay }
        elseif ($hl.Range.Row -eq $row2) { $text2 = $hl.TextToDisplay }
    }
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Row -eq $row1) { $hl.TextToDisplay = $text2 }
        elseif ($hl.Range.Row -eq $row2) { $hl.TextToDisplay = $text1 }
    }
}

$wb = $excel.ActiveWorkbook

# --- Overview sheet: File Name (A), Path And Name (B, hyperlinked),
#     zh-cn (E), de-de (F), Latest HO Xliff Generate Date (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
Swap-RowCells $wsOverview 3 4 @("A","B","E","F","G")
Swap-HyperlinkDisplay $wsOverview 3 4

# --- zh-cn sheet: Source File Name (A, hyperlinked), Status (C),
#     Latest Handoff File (G), Latest Handoff Datetime (H) ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
Swap-RowCells $wsZhCn 3 4 @("A","C","G","H")
Swap-HyperlinkDisplay $wsZhCn 3 4

# --- de-de sheet: same layout as zh-cn ---
$wsDeDe = $wb.Worksheets.Item("de-de")
Swap-RowCells $wsDeDe 3 4 @("A","C","G","H")
Swap-HyperlinkDisplay $wsDeDe 3 4
